$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "quartile 4" column from the Developed group (old column F) --
# shifts old G:K left into F:J.
$ws.Columns("F").Delete()

# Drop the "quartile 4" column from the Emerging group (now column J,
# formerly K) -- shifts nothing further since it was the last column.
$ws.Columns("J").Delete()

# The forward-return row now holds freshly processed values rather than a
# simple shift of the old ones, so set them explicitly.
$ws.Range("B4").Value = 0.01030003143891282
$ws.Range("C4").Value = 0.00700879275723626
$ws.Range("D4").Value = 0.006289948395120337
$ws.Range("E4").Value = 0.006867600030428913
$ws.Range("F4").Value = 0.0176036064588034
$ws.Range("G4").Value = 0.01169426275712809
$ws.Range("H4").Value = 0.01053300212894502
$ws.Range("I4").Value = 0.007510228626348137

Write-Output "done"
